# InventarioDatos-ScoringPJ.xlsx style edit
# 1) Remove the redundant duplicate "Fuente de Datos" column (old column C)
# 2) Populate the new columns (Proposito .. Derechos y Restricciones) for rows 2-6
# 3) Fix the wording of the "Descripcion" text for rows 4-6
# 4) Add a brand new data row (row 7, DANE source)
# 5) Re-point the hyperlinks at the Url column (now column C)
# 6) Left-align every data/header cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Drop the duplicate column (old "C", a copy of "Fuente de Datos")
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).Delete()

# ---------------------------------------------------------------------------
# 2) Fix wording for the "Descripcion" column (D) on rows 4-6
# ---------------------------------------------------------------------------
$nuevaDescripcionBancos = "Esta información se usará para simular el comportamiento con los bancos del sitema financiero (diferentes a los del Banco de Bogotá) de los clientes del perfil a impactar."
$ws.Range("D4").Value = $nuevaDescripcionBancos
$ws.Range("D5").Value = $nuevaDescripcionBancos
$ws.Range("D6").Value = $nuevaDescripcionBancos

# ---------------------------------------------------------------------------
# 3) Fill in the new columns E (Proposito) .. M (Derechos y Restricciones)
#    for the five existing data rows (2-6)
# ---------------------------------------------------------------------------

# Row 2 - Superintendencia de Sociedades (NIIF)
$ws.Range("E2").Value = "Segmentar la información financiera del cliente en indicadores, promedios, etc."
$ws.Range("F2").Value = "Superintendencia de Sociedades"
$ws.Range("G2").Value = "Superintendencia de Sociedades"
$ws.Range("H2").Value = "Conocer la situación financiera actual y pasada del cliente"
$ws.Range("I2").Value = 2021
$ws.Range("J2").Value = "Anual"
$ws.Range("K2").Value = "Bases de Datos"
$ws.Range("M2").Value = "Información Pública"

# Row 3 - Superintendencia de Sociedades (SIREM)
$ws.Range("E3").Value = "Segmentar la información financiera del cliente en indicadores, promedios, etc."
$ws.Range("F3").Value = "Superintendencia de Sociedades"
$ws.Range("G3").Value = "Superintendencia de Sociedades"
$ws.Range("H3").Value = "Conocer la situación financiera actual y pasada del cliente"
$ws.Range("I3").Value = 2021
$ws.Range("J3").Value = "Anual"
$ws.Range("K3").Value = "Bases de Datos"
$ws.Range("M3").Value = "Información Pública"

# Row 4 - Datacredito
$ws.Range("E4").Value = 50
$ws.Range("F4").Value = "Datacrédito"
$ws.Range("G4").Value = "Datacrédito"
$ws.Range("H4").Value = "Ver el comportamiento financiero del grupo objetivo"
$ws.Range("I4").Value = "N/A"
$ws.Range("J4").Value = "Según consulta"
$ws.Range("K4").Value = "Bases de Datos"
$ws.Range("M4").Value = "Información Privadas de pago por consulta"

# Row 5 - CIFIN
$ws.Range("E5").Value = 50
$ws.Range("F5").Value = "CIFIN"
$ws.Range("G5").Value = "CIFIN"
$ws.Range("H5").Value = "Ver el comportamiento financiero del grupo objetivo"
$ws.Range("I5").Value = "N/A"
$ws.Range("J5").Value = "Según consulta"
$ws.Range("K5").Value = "Bases de Datos"
$ws.Range("M5").Value = "Información Privadas de pago por consulta"

# Row 6 - Banco de Bogota (vector moras)
$ws.Range("E6").Value = 50
$ws.Range("F6").Value = "M.I.S Banco de Bogotá"
$ws.Range("G6").Value = "Banco de Bogotá"
$ws.Range("H6").Value = "Ver el comportamiento del grupo objetivo con la entidad objetivo"
$ws.Range("I6").Value = "N/A"
$ws.Range("J6").Value = "Mensual"
$ws.Range("K6").Value = "Bases de Datos"
$ws.Range("M6").Value = "Imformación privada sin derech de distribución"

# Row 6 originally only had the borderless "F6:N6" template style; F6/G6 keep
# that borderless look (as in the final file) but the rest of the newly
# filled-in cells need the normal bordered/wrapped look used elsewhere, so
# copy that formatting over from the row above (values are left untouched).
$ws.Range("H5").Copy()
$ws.Range("E6").PasteSpecial(-4122)
$ws.Range("H6").PasteSpecial(-4122)
$ws.Range("I6").PasteSpecial(-4122)
$ws.Range("J6").PasteSpecial(-4122)
$ws.Range("K6").PasteSpecial(-4122)
$ws.Range("M6").PasteSpecial(-4122)

Write-Host "stage1 ok"
